$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    "B2" = 0.2230769230769231
    "C2" = 0.5230769230769231
    "J2" = 0.003846153846153846
    "P2" = 0.15
    "S2" = 0.1
    "B3" = 0.007246376811594203
    "C3" = 0.01449275362318841
    "J3" = 0.01449275362318841
    "P3" = 0.7681159420289855
    "S3" = 0.1956521739130435
    "P4" = 0.6904761904761905
    "S4" = 0.3095238095238095
    "B6" = 0.1183431952662722
    "D6" = 0.01775147928994083
    "F6" = 0.0650887573964497
    "J6" = 0.1893491124260355
    "O6" = 0.02958579881656805
    "Q6" = 0.1479289940828402
    "R6" = 0.04733727810650887
    "S6" = 0.3846153846153846
    "B7" = 0.0898876404494382
    "D7" = 0.01685393258426966
    "F7" = 0.05617977528089887
    "J7" = 0.1179775280898876
    "O7" = 0.03932584269662921
    "Q7" = 0.151685393258427
    "R7" = 0.07303370786516854
    "S7" = 0.4550561797752809
    "B8" = 0.1088082901554404
    "D8" = 0.01813471502590673
    "F8" = 0.04404145077720207
    "J8" = 0.1217616580310881
    "O8" = 0.0310880829015544
    "Q8" = 0.1580310880829016
    "R8" = 0.07772020725388601
    "S8" = 0.4404145077720207
    "B9" = 0.05263157894736842
    "D9" = 0.03289473684210526
    "F9" = 0.07894736842105263
    "J9" = 0.1513157894736842
    "O9" = 0.006578947368421052
    "Q9" = 0.2171052631578947
    "R9" = 0.03289473684210526
    "S9" = 0.4276315789473684
    "B10" = 0.1114341085271318
    "D10" = 0.0251937984496124
    "F10" = 0.07170542635658915
    "J10" = 0.1143410852713178
    "O10" = 0.01937984496124031
    "Q10" = 0.2354651162790698
    "R10" = 0.0562015503875969
    "S10" = 0.3662790697674418
    "G11" = 0.1532567049808429
    "J11" = 0.06130268199233716
    "K11" = 0.2030651340996169
    "L11" = 0.5747126436781609
    "S11" = 0.007662835249042145
    "G12" = 0.7388535031847133
    "J12" = 0.1592356687898089
    "K12" = 0.02547770700636943
    "L12" = 0.03821656050955414
    "S12" = 0.03821656050955414
    "G13" = 0.6571428571428571
    "J13" = 0.2857142857142857
    "S13" = 0.05714285714285714
    "F15" = 0.015625
    "H15" = 0.1145833333333333
    "I15" = 0.08333333333333333
    "J15" = 0.375
    "K15" = 0.08333333333333333
    "M15" = 0.01041666666666667
    "S15" = 0.2552083333333333
    "F16" = 0.005952380952380952
    "H16" = 0.1428571428571428
    "I16" = 0.07142857142857142
    "J16" = 0.4285714285714285
    "K16" = 0.1369047619047619
    "M16" = 0.01785714285714286
    "O16" = 0.05357142857142857
    "S16" = 0.1428571428571428
    "F17" = 0.01272264631043257
    "H17" = 0.2061068702290076
    "I17" = 0.08651399491094147
    "J17" = 0.4071246819338422
    "K17" = 0.1195928753180662
    "M17" = 0.01272264631043257
    "N17" = 0.002544529262086514
    "O17" = 0.04325699745547074
    "S17" = 0.1094147582697201
    "H18" = 0.1565217391304348
    "I18" = 0.06956521739130435
    "J18" = 0.4
    "K18" = 0.1304347826086956
    "M18" = 0.01739130434782609
    "N18" = 0.008695652173913044
    "O18" = 0.06956521739130435
    "S18" = 0.1478260869565217
    "F19" = 0.01197053406998158
    "H19" = 0.2274401473296501
    "I19" = 0.07642725598526703
    "J19" = 0.3710865561694291
    "K19" = 0.09208103130755065
    "M19" = 0.02394106813996317
    "N19" = 0.0009208103130755065
    "O19" = 0.0718232044198895
    "S19" = 0.1243093922651934
}

foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}
